$d = $word.ActiveDocument

# Paragraph 2: "Test link before bookmark : ..." -- update rsidR on the
# field-code runs (begin/instrText/separate/result/end) that wrap the
# "a reference to bookmark1" REF field result.
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
if ($r2.Text -notlike "Test link before bookmark*") {
    throw "Paragraph 2 did not match expected content: $($r2.Text)"
}
$p2xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidP="009168BC" w:rsidR="00E02A2B" w:rsidRDefault="00E02A2B" w:rsidRPr="007F2DB9"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Test link before bookmark : </w:t></w:r><w:r><w:rPr><w:b w:val="true"/><w:color w:val="FF0000"/></w:rPr><w:t>dangling reference for bookmark bookmark1</w:t></w:r><w:r/><w:r w:rsidR="024413B25F474757AF555E68CE726323"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="024413B25F474757AF555E68CE726323"><w:instrText xml:space="preserve"/></w:r><w:r w:rsidR="024413B25F474757AF555E68CE726323"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="024413B25F474757AF555E68CE726323"><w:rPr><w:b w:val="true"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r><w:r w:rsidR="024413B25F474757AF555E68CE726323"><w:fldChar w:fldCharType="end"/></w:r></w:p>
'@
$r2.InsertXML($p2xml)

# Paragraph 3: "Test bookmark : divOp(...) failed: / by zero" -- append the
# full Java stack trace to the error message run.
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
if ($r3.Text -notlike "Test*bookmark*divOp*") {
    throw "Paragraph 3 did not match expected content: $($r3.Text)"
}
$p3xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidP="009168BC" w:rsidR="00C52979" w:rsidRDefault="00E02A2B" w:rsidRPr="007F2DB9"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Test</w:t></w:r><w:r w:rsidR="00C52979" w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>bookmark</w:t></w:r><w:r w:rsidR="00C52979" w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> : </w:t></w:r><w:r w:rsidR="00C52979"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="002D1E44" w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:instrText>m</w:instrText></w:r><w:r w:rsidR="002848A5" w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:instrText>:</w:instrText></w:r><w:r w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:instrText>bookmark</w:instrText></w:r><w:r w:rsidR="003D27D6" w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:instrText>1/0</w:instrText></w:r><w:r w:rsidR="00C52979"><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:b w:val="true"/><w:color w:val="FF0000"/></w:rPr><w:t>divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:
	/ by zero
java.lang.ArithmeticException: / by zero
	at org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)
	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)
	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)
	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)
	at java.lang.reflect.Method.invoke(Method.java:498)
	at org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)
	at org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)
	at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)
	at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)
	at org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)
	at org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)
	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
	at org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)
	at org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBookmark(M2DocEvaluator.java:1168)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBookmark(M2DocEvaluator.java:1)
	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:253)
	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1034)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)
	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:183)
	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:297)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:1)
	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:201)
	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:259)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)
	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:246)
	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:252)
	at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:691)
	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:396)
	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:318)
	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)
	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)
	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)
	at java.lang.reflect.Method.invoke(Method.java:498)
	at org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)
	at org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)
	at org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)
	at org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)
	at org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)
	at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)
	at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)
	at org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)
	at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)
	at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)
	at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)
	at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)
	at org.junit.runners.ParentRunner.run(ParentRunner.java:363)
	at org.junit.runners.Suite.runChild(Suite.java:128)
	at org.junit.runners.Suite.runChild(Suite.java:27)
	at org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)
	at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)
	at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)
	at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)
	at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)
	at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)
	at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)
	at org.junit.runners.ParentRunner.run(ParentRunner.java:363)
	at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)
	at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:459)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:675)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:382)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:192)
</w:t><w:br/></w:r><w:r w:rsidR="003D27D6"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="003D27D6" w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="002D1E44" w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:instrText>m</w:instrText></w:r><w:r w:rsidR="003D27D6" w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:instrText>:end</w:instrText></w:r><w:r w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:instrText>bookmark</w:instrText></w:r><w:r w:rsidR="003D27D6" w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="003D27D6"><w:fldChar w:fldCharType="end"/></w:r></w:p>
'@
$r3.InsertXML($p3xml)

# Paragraph 4: "Test link after bookmark : ..." -- update rsidR on the
# field-code runs the same way as paragraph 2.
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
if ($r4.Text -notlike "Test link after bookmark*") {
    throw "Paragraph 4 did not match expected content: $($r4.Text)"
}
$p4xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidP="00E02A2B" w:rsidR="00E02A2B" w:rsidRDefault="00E02A2B" w:rsidRPr="007F2DB9"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Test link after bookmark : </w:t></w:r><w:r><w:rPr><w:b w:val="true"/><w:color w:val="FF0000"/></w:rPr><w:t>dangling reference for bookmark bookmark1</w:t></w:r><w:r/><w:r w:rsidR="024413B25F474757AF555E68CE726323"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="024413B25F474757AF555E68CE726323"><w:instrText xml:space="preserve"/></w:r><w:r w:rsidR="024413B25F474757AF555E68CE726323"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="024413B25F474757AF555E68CE726323"><w:rPr><w:b w:val="true"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r><w:r w:rsidR="024413B25F474757AF555E68CE726323"><w:fldChar w:fldCharType="end"/></w:r><w:r w:rsidR="00D0546C" w:rsidRPr="007F2DB9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@
$r4.InsertXML($p4xml)

Write-Host "Done."
